$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General")

# Insert a new row 4 (pushes old rows 4-6 down to 5-7) and populate it
$ws.Rows.Item(4).Insert()

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Insecure Web Administration"

$textSummary = @'
Most devices offer some form of web application to provision and administer the device. These interfaces are vulnerable to the same risks as enterprise applications or Internet sites. However, it is true that IoT admin interfaces have not been generally subjected to the same security assessment as those more familiar targets.
In fact, it is true that these admin interfaces are at increased risk of attack because the owner of the device has physical access. As such the firmware can be retrieved which includes the code for the web application. Attackers have time to sift through this code to find weaknesses while they do not have this luxury for traditional web sites in most cases.
Common web application flaws in embedded admin interfaces include:
• Remote code execution
• Local/Remote File Includes
• Directory Traversal
• Cross Site Scripting (XSS)
• Denial of Service
Every application level flaw is applicable but these are likely within IoT admin interfaces.

'@
$ws.Range("C4").Value = $textSummary

$textRecommendation = @'
Educate developers on the risks of common web application flaws. Security engaged engineers are effective to limit risks.
Leverage static code review tools which may help you locate flaws early.
However, it is also key to have your administration applications assessed by a trusted 3rd party. This should not be an automated vulnerability assessment it needs to include manual assessment to find the best results.
To make this process as cost effective as possible we would recommend providing the source code for the interface along with access to the underlying operating system. 
To ensure coverage it is best to start from a position of knowing every file accessible via the web server. Thinking that someone will “never find” content that is not linked to directly is wrong in this case. Since the firmware will be accessible the attacker will have access anyway if they look in the right places.
While classic web application testing is the skillset required. For the purposes of embedded devices it is best to provide as much information as possible so the audit can be comprehensive.

'@
$ws.Range("D4").Value = $textRecommendation

$ws.Rows.Item(4).RowHeight = 289.2

# Make "General" the active sheet/tab, select E4
$ws.Activate() | Out-Null
$ws.Range("E4").Select() | Out-Null

